$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (new report week / issue number) ---
$ws.Range("A8").Value = "Volume 31   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/5/2024  Through  8/11/2024"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("F14").Value = "0"
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = -100
$ws.Range("I14").Value = 6
$ws.Range("J14").Value = 7
$ws.Range("K14").Value = -14.285714285714
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = -82.857142857142
$ws.Range("G15").Value = "0"
$ws.Range("H15").Value = "***.*"
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -82.142857142857
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -14.285714285714
$ws.Range("F16").Value = 25
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 38.888888888888
$ws.Range("I16").Value = 144
$ws.Range("J16").Value = 127
$ws.Range("K16").Value = 13.385826771653
$ws.Range("L16").Value = 7.462686567164
$ws.Range("M16").Value = -19.553072625698
$ws.Range("N16").Value = -73.480662983425
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = -17.777777777777
$ws.Range("I17").Value = 301
$ws.Range("J17").Value = 290
$ws.Range("K17").Value = 3.793103448275
$ws.Range("L17").Value = -9.609609609609
$ws.Range("M17").Value = 49.751243781094
$ws.Range("N17").Value = -51.056910569105
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -47.368421052631
$ws.Range("I18").Value = 78
$ws.Range("J18").Value = 114
$ws.Range("K18").Value = -31.578947368421
$ws.Range("L18").Value = -37.096774193548
$ws.Range("M18").Value = 8.333333333333
$ws.Range("N18").Value = -86.912751677852
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = -45.652173913043
$ws.Range("I19").Value = 212
$ws.Range("J19").Value = 261
$ws.Range("K19").Value = -18.773946360153
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 14.594594594594
$ws.Range("N19").Value = -8.225108225108
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -71.428571428571
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -60
$ws.Range("I20").Value = 57
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = -18.571428571428
$ws.Range("L20").Value = 11.764705882352
$ws.Range("M20").Value = 83.870967741935
$ws.Range("N20").Value = -65.868263473053
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -42.105263157894
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = 145
$ws.Range("H21").Value = -28.965517241379
$ws.Range("I21").Value = 808
$ws.Range("J21").Value = 887
$ws.Range("K21").Value = -8.90642615558
$ws.Range("L21").Value = -7.867730900798
$ws.Range("M21").Value = 17.101449275362
$ws.Range("N21").Value = -63.976816763263
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = -36.363636363636
$ws.Range("L22").Value = -50
$ws.Range("M22").Value = 40
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -40
$ws.Range("F23").Value = 20
$ws.Range("G23").Value = 23
$ws.Range("H23").Value = -13.043478260869
$ws.Range("I23").Value = 172
$ws.Range("J23").Value = 145
$ws.Range("K23").Value = 18.620689655172
$ws.Range("L23").Value = 25.547445255474
$ws.Range("M23").Value = 70.29702970297
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -45.833333333333
$ws.Range("F24").Value = 45
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = -48.275862068965
$ws.Range("I24").Value = 513
$ws.Range("J24").Value = 535
$ws.Range("K24").Value = -4.11214953271
$ws.Range("L24").Value = -3.571428571428
$ws.Range("M24").Value = 33.942558746736
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = -9.090909090909
$ws.Range("I25").Value = 98
$ws.Range("J25").Value = 107
$ws.Range("K25").Value = -8.411214953271
$ws.Range("L25").Value = -30.496453900709
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 300
$ws.Range("F26").Value = 76
$ws.Range("G26").Value = 36
$ws.Range("H26").Value = 111.111111111111
$ws.Range("I26").Value = 443
$ws.Range("J26").Value = 361
$ws.Range("K26").Value = 22.714681440443
$ws.Range("L26").Value = 24.788732394366
$ws.Range("M26").Value = -19.746376811594
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 17
$ws.Range("K27").Value = -45.16129032258
$ws.Range("L27").Value = -41.379310344827
$ws.Range("C28").Value = "0"
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -57.142857142857
$ws.Range("J28").Value = 35
$ws.Range("K28").Value = -14.285714285714
$ws.Range("L28").Value = -43.396226415094
$ws.Range("D29").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = -80
$ws.Range("J29").Value = 22
$ws.Range("K29").Value = -31.818181818181
$ws.Range("M29").Value = -42.307692307692
$ws.Range("N29").Value = -80.76923076923
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = -75
$ws.Range("J30").Value = 19
$ws.Range("K30").Value = -36.842105263157
$ws.Range("M30").Value = -42.857142857142
$ws.Range("N30").Value = -83.333333333333
$ws.Range("F33").Value = "0"
